$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data appended below the existing data (rows 5-7).
# Column A: Date (serial date-time values, formatted like existing rows)
# Column B: Method ("Random")
# Columns C-M: numeric metrics

$rows = @(
    @{ Row = 5; Date = 42602.514328703706; H = 48; I = 52; L = 32; M = 68 },
    @{ Row = 6; Date = 42602.516215277778; H = 20; I = 80; L = 75; M = 25 },
    @{ Row = 7; Date = 42602.517106481479; H = 30; I = 70; L = 67; M = 33 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = "Random"
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
    $ws.Cells.Item($row, 5).Value = 0
    $ws.Cells.Item($row, 6).Value = 0
    $ws.Cells.Item($row, 7).Value = 0
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}
